$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52 (shifts existing rows 52-93 down to 53-94)
$ws.Rows(52).Insert()

# Populate the new row 52 with the latest weekly price entry.
# Columns A,B,C,E,F,G,H,I,Q,R mirror the rest of the series (unchanged metadata);
# D (Fecha), J (Volumen), K/L/M (precios), N (unidad), O (origen) and P (precio $/Kg)
# carry the new week's values.
$ws.Range("A52").Value = 4
$ws.Range("B52").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C52").Value = "Los Lagos"
$ws.Range("D52").Value = 44789
$ws.Range("E52").Value = 10
$ws.Range("F52").Value = 100112031
$ws.Range("G52").Value = "Poroto verde"
$ws.Range("H52").Value = "Magnum"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 40
$ws.Range("K52").Value = 37000
$ws.Range("L52").Value = 37000
$ws.Range("M52").Value = 37000
$ws.Range("N52").Value = "`$/malla 25 kilos"
$ws.Range("O52").Value = "Perú"
$ws.Range("P52").Value = 1480
$ws.Range("Q52").Value = 25
$ws.Range("R52").Value = "Hortaliza"
